# Adds a new "Commit 5" block (rows 87-102) to Sheet1, mirroring the
# existing "Commit 4" block (rows 70-85), per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Values & formulas ------------------------------------------------

# Section title
$ws.Range("A87").Value = "Commit 5"

# "MARS Tool Output" / "Calulations" header band
$ws.Range("A88").Value = "MARS Tool Output"
$ws.Range("B88").Value = ""
$ws.Range("C88").Value = ""
$ws.Range("D88").Value = "Calulations"

# Instruction Statistics Tool
$ws.Range("A90").Value = "Instruction Statistics Tool"

$ws.Range("A91").Value = "Instruction type"
$ws.Range("B91").Value = "Count"
$ws.Range("D91").Value = "Adjusted count"
$ws.Range("E91").Value = "CPI"
$ws.Range("F91").Value = "Total cycles"

$ws.Range("A92").Value = "ALU"
$ws.Range("B92").Value = 3608
$ws.Range("D92").Formula = "=B92"
$ws.Range("E92").Value = 1
$ws.Range("F92").Formula = "=D92*E92"

$ws.Range("A93").Value = "Jump"
$ws.Range("B93").Value = 72
$ws.Range("D93").Formula = "=B93"
$ws.Range("E93").Value = 1
$ws.Range("F93").Formula = "=D93*E93"

$ws.Range("A94").Value = "Branch"
$ws.Range("B94").Value = 915
$ws.Range("D94").Formula = "=B94"
$ws.Range("E94").Value = 2
$ws.Range("F94").Formula = "=D94*E94"

$ws.Range("A95").Value = "Memory"
$ws.Range("B95").Value = 617

$ws.Range("A96").Value = "Other"
$ws.Range("B96").Value = 712
$ws.Range("D96").Formula = "=B96-(B100+B101-B95)"
$ws.Range("E96").Value = 5
$ws.Range("F96").Formula = "=D96*E96"

# Data Cache Simulation Tool
$ws.Range("A98").Value = "Data Cache Simulation Tool"

$ws.Range("A99").Value = "Access"
$ws.Range("B99").Value = "Count"

$ws.Range("A100").Value = "Cache hit"
$ws.Range("B100").Value = 359
$ws.Range("D100").Formula = "=B100"
$ws.Range("E100").Value = 2
$ws.Range("F100").Formula = "=D100*E100"

$ws.Range("A101").Value = "Cache miss"
$ws.Range("B101").Value = 350
$ws.Range("D101").Formula = "=B101"
$ws.Range("E101").Value = 40
$ws.Range("F101").Formula = "=D101*E101"

$ws.Range("F102").Formula = "=SUM(F92:F101)"

# ---- Formatting (reuse the styles from the "Commit 4" block above) ---

function Copy-Format([string]$src, [string]$dst) {
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial(-4122) | Out-Null
}

Copy-Format "A71" "A88"
Copy-Format "B71" "B88"
Copy-Format "C71" "C88"
Copy-Format "D71" "D88"

Copy-Format "A73" "A90"

Copy-Format "A74" "A91"
Copy-Format "B74" "B91"
Copy-Format "D74" "D91"
Copy-Format "E74" "E91"
Copy-Format "F74" "F91"

Copy-Format "A81" "A98"

Copy-Format "A82" "A99"
Copy-Format "B82" "B99"

Copy-Format "F85" "F102"

$excel.CutCopyMode = $false

# ---- Selection / view state --------------------------------------------

$ws.Range("B102").Select() | Out-Null
